$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append a new data row (row 96) to the daily log sheet.
# Column A holds a date-like text label, so force text formatting first
# to avoid Excel auto-converting the string into a date serial number,
# then reset the cell style back to Normal so no extra formatting is
# introduced (matching the plain, unstyled data rows above it).
$row = 96

$cellA = $ws.Cells.Item($row, 1)
$cellA.NumberFormat = "@"
$cellA.Value = "2025/10/13"
$cellA.Style = "Normal"

$ws.Cells.Item($row, 2).Value = "月"
$ws.Cells.Item($row, 3).Value = 0
$ws.Cells.Item($row, 4).Value = 44
